# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp (shared string reused at A1)
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 08:04"

# --- Simple numeric corrections (no row reordering) ---

# Alemania (row 10)
$ws.Range("D10").Value = 143300
$ws.Range("E10").Value = 19778

# India (row 17)
$ws.Range("B17").Value = 59765
$ws.Range("C17").Value = 70
$ws.Range("D17").Value = 17897
$ws.Range("E17").Value = 39882
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1986

# Hungria (row 66)
$ws.Range("B66").Value = 3213
$ws.Range("C66").Value = 35
$ws.Range("D66").Value = 904
$ws.Range("E66").Value = 1904
$ws.Range("G66").Value = 13
$ws.Range("H66").Value = 405

# Uzbekistan (row 72)
$ws.Range("B72").Value = 2336
$ws.Range("C72").Value = 11
$ws.Range("E72").Value = 551

# Bulgaria (row 79)
$ws.Range("B79").Value = 1911
$ws.Range("C79").Value = 39
$ws.Range("D79").Value = 422
$ws.Range("E79").Value = 1401
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 88

# --- Re-order Kirguistan ahead of Somalia / Letonia (rows 96-98) ---
# Row 96 becomes Kirguistan with freshly updated figures
$ws.Range("A96").Value = "Kirguistan"
$ws.Range("B96").Value = 931
$ws.Range("C96").Value = 25
$ws.Range("D96").Value = 658
$ws.Range("E96").Value = 261
$ws.Range("F96").Value = 13
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 12

# Row 97 becomes Somalia (old row 96 figures, unchanged)
$ws.Range("A97").Value = "Somalia"
$ws.Range("B97").Value = 928
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 106
$ws.Range("E97").Value = 778
$ws.Range("F97").Value = 2
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 44

# Row 98 becomes Letonia (old row 97 figures, unchanged)
$ws.Range("A98").Value = "Letonia"
$ws.Range("B98").Value = 928
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 464
$ws.Range("E98").Value = 446
$ws.Range("F98").Value = 2
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 18

# --- Re-order Belice ahead of Nueva Caledonia (rows 192-193) ---
# Row 192 becomes Belice (old row 193 figures)
$ws.Range("A192").Value = "Belice"
$ws.Range("B192").Value = 18
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 16
$ws.Range("E192").Value = 0
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 2

# Row 193 becomes Nueva Caledonia (old row 192 figures)
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("B193").Value = 18
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 18
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0
